# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the zh-cn /
# de-de handoff packages have moved from "In Translation" to
# "Ready for handoff", refreshes the associated handoff timestamps, and
# widens the status/datetime columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$colWidth  = 16.333333333333336   # closest reachable width to the new 17.22-char target

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) + HO xliff date (G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-17 12:36:55"

$wsOverview.Columns.Item(5).ColumnWidth = $colWidth
$wsOverview.Columns.Item(6).ColumnWidth = $colWidth

# ---------------------------------------------------------------------
# zh-cn sheet: Status (C) + Latest Handoff Datetime (H)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-17 12:36:50"

$wsZhCn.Columns.Item(3).ColumnWidth = $colWidth

# ---------------------------------------------------------------------
# de-de sheet: Status (C) + Latest Handoff Datetime (H)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-17 12:36:55"

$wsDeDe.Columns.Item(3).ColumnWidth = $colWidth
